$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for first row updated
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 15:07:37"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime for first row updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 15:07:32"
$wsZhCn.Range("K2").Value = "2016-08-16 15:07:49"

# Sheet "de-de": Correspond Handback DateTime for first row updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-16 15:07:56"
